{"js": "// Replace the date line and each \"AxB=\" expression in the practice table\n// with its updated counterpart. Each \"before\" string is unique in the\n// document, so a simple ordered search-and-replace pass is unambiguous.\nconst replacements = [\n  [\"2025-10-21 Tuesday\", \"2025-10-22 Wednesday\"],\n  [\"370\u00d77=\", \"650\u00d75=\"],\n  [\"634\u00d75=\", \"279\u00d78=\"],\n  [\"407\u00d76=\", \"566\u00d78=\"],\n  [\"424\u00d76=\", \"756\u00d72=\"],\n  [\"178\u00d78=\", \"664\u00d75=\"],\n  [\"675\u00d72=\", \"444\u00d76=\"],\n  [\"554\u00d72=\", \"898\u00d74=\"],\n  [\"374\u00d77=\", \"940\u00d77=\"],\n  [\"313\u00d76=\", \"297\u00d74=\"],\n  [\"904\u00d73=\", \"258\u00d79=\"],\n  [\"687\u00d73=\", \"118\u00d73=\"],\n  [\"919\u00d78=\", \"788\u00d79=\"],\n  [\"446\u00d76=\", \"458\u00d75=\"],\n  [\"225\u00d79=\", \"898\u00d78=\"],\n  [\"879\u00d78=\", \"583\u00d79=\"],\n  [\"113\u00d77=\", \"834\u00d73=\"],\n  [\"634\u00d72=\", \"577\u00d77=\"],\n  [\"499\u00d78=\", \"606\u00d78=\"],\n  [\"249\u00d74=\", \"854\u00d74=\"],\n  [\"384\u00d79=\", \"523\u00d78=\"],\n  [\"233\u00d79=\", \"323\u00d76=\"],\n  [\"714\u00d78=\", \"822\u00d74=\"],\n  [\"941\u00d72=\", \"930\u00d72=\"],\n  [\"431\u00d79=\", \"387\u00d75=\"],\n  [\"707\u00d72=\", \"842\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"AxB=\" expression in the practice table\n# with its updated counterpart. Each \"before\" string is unique in the\n# document, so a simple ordered Find/Replace pass is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-21 Tuesday\", \"2025-10-22 Wednesday\"),\n    @(\"370\u00d77=\", \"650\u00d75=\"),\n    @(\"634\u00d75=\", \"279\u00d78=\"),\n    @(\"407\u00d76=\", \"566\u00d78=\"),\n    @(\"424\u00d76=\", \"756\u00d72=\"),\n    @(\"178\u00d78=\", \"664\u00d75=\"),\n    @(\"675\u00d72=\", \"444\u00d76=\"),\n    @(\"554\u00d72=\", \"898\u00d74=\"),\n    @(\"374\u00d77=\", \"940\u00d77=\"),\n    @(\"313\u00d76=\", \"297\u00d74=\"),\n    @(\"904\u00d73=\", \"258\u00d79=\"),\n    @(\"687\u00d73=\", \"118\u00d73=\"),\n    @(\"919\u00d78=\", \"788\u00d79=\"),\n    @(\"446\u00d76=\", \"458\u00d75=\"),\n    @(\"225\u00d79=\", \"898\u00d78=\"),\n    @(\"879\u00d78=\", \"583\u00d79=\"),\n    @(\"113\u00d77=\", \"834\u00d73=\"),\n    @(\"634\u00d72=\", \"577\u00d77=\"),\n    @(\"499\u00d78=\", \"606\u00d78=\"),\n    @(\"249\u00d74=\", \"854\u00d74=\"),\n    @(\"384\u00d79=\", \"523\u00d78=\"),\n    @(\"233\u00d79=\", \"323\u00d76=\"),\n    @(\"714\u00d78=\", \"822\u00d74=\"),\n    @(\"941\u00d72=\", \"930\u00d72=\"),\n    @(\"431\u00d79=\", \"387\u00d75=\"),\n    @(\"707\u00d72=\", \"842\u00d73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
